$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

# ---------------------------------------------------------------------------
# SPRINT 5 section (new rows 92, 94-104) on the "Sprint" sheet, mirroring the
# layout already used for SPRINT 1-4 further up the sheet.
# ---------------------------------------------------------------------------

# Section title
$ws.Range("A92").Value = "SPRINT 5"

# Column headers (row 94)
$ws.Range("B94").Value = "Task (Definition of Done)"
$ws.Range("C94").Value = "Waktu perkiraan (baseline)"
$ws.Range("D94").Value = "Hari 1"
$ws.Range("E94").Value = "Hari 2"
$ws.Range("F94").Value = "Hari 3"
$ws.Range("G94").Value = "Total"

# Task rows 95-100
$ws.Range("B95").Value = "Mendesign dan menuliskan pengertian Bus SCSI pada slide"
$ws.Range("C95").Value = 1
$ws.Range("D95").Value = 1
$ws.Range("E95").Value = 0
$ws.Range("F95").Value = 0
$ws.Range("G95").Formula = "=SUM(D95:F95)"

$ws.Range("B96").Value = "Membuat slide baru dengan menambahkan gambar Bus SCSI didalamnya"
$ws.Range("C96").Value = 1
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0

$ws.Range("B97").Value = "Mendesign dan menuliskan pengertian Bus Prosessor pada slide"
$ws.Range("C97").Value = 1
$ws.Range("D97").Value = 1
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0

$ws.Range("B98").Value = "Membuat slide baru dengan menambahkan gambar Bus Prosessor didalamnya"
$ws.Range("C98").Value = 1
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0

$ws.Range("B99").Value = "Mendesign dan menuliskan pengertian Bus AGP pada slide"
$ws.Range("C99").Value = 1
$ws.Range("D99").Value = 1
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0

$ws.Range("B100").Value = "Membuat slide baru dengan menambahkan gambar Bus AGP didalamnya"
$ws.Range("C100").Value = 1
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 0
$ws.Range("F100").Value = 0

# Shared "Total" formula fill-down for G96:G100 (mirrors G95's formula)
$ws.Range("G96:G100").Formula = "=SUM(D96:F96)"

# Blank spacer row 101 — clear any stray formatting/content (nothing to set).

# Secondary header row (102) for the summary block
$ws.Range("C102").Value = "Hari 0"
$ws.Range("D102").Value = "Hari 1"
$ws.Range("E102").Value = "Hari 2"
$ws.Range("F102").Value = "Hari 3"

# Summary rows 103 (actual time) and 104 (baseline estimate)
$ws.Range("B103").Value = "Waktu sebenarnya"
$ws.Range("C103").Formula = "=SUM(C95:C101)"
$ws.Range("D103").Formula = "=C103-(SUM(D95:D101))"
$ws.Range("E103").Formula = "=D103-(SUM(E95:E101))"
$ws.Range("F103").Formula = "=E103-(SUM(F95:F101))"

$ws.Range("B104").Value = "Waktu perkiraan (baseline)"
$ws.Range("C104").Formula = "=SUM(C95:C101)"
$ws.Range("D104").Formula = "=C104-(C104/3)"
$ws.Range("E104").Formula = "=D104-(C104/3)"
$ws.Range("F104").Formula = "=E104-(C104/3)"

# ---------------------------------------------------------------------------
# Formatting — copy the styling already used for the SPRINT 4 block
# (rows 69-84) onto the new SPRINT 5 block (rows 92-104) so it matches the
# rest of the sheet.
# ---------------------------------------------------------------------------

$ws.Range("A69").Copy() | Out-Null
$ws.Range("A92").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B71:G71").Copy() | Out-Null
$ws.Range("B94:G94").PasteSpecial(-4122) | Out-Null

$ws.Range("B72:G72").Copy() | Out-Null
$ws.Range("B95:G95").PasteSpecial(-4122) | Out-Null
$ws.Range("B96:G96").PasteSpecial(-4122) | Out-Null
$ws.Range("B97:G97").PasteSpecial(-4122) | Out-Null
$ws.Range("B98:G98").PasteSpecial(-4122) | Out-Null
$ws.Range("B99:G99").PasteSpecial(-4122) | Out-Null
$ws.Range("B100:G100").PasteSpecial(-4122) | Out-Null

$ws.Range("B81:G81").Copy() | Out-Null
$ws.Range("B101:G101").PasteSpecial(-4122) | Out-Null

$ws.Range("C82:G82").Copy() | Out-Null
$ws.Range("C102:G102").PasteSpecial(-4122) | Out-Null

$ws.Range("B83:G83").Copy() | Out-Null
$ws.Range("B103:G103").PasteSpecial(-4122) | Out-Null

$ws.Range("B84:F84").Copy() | Out-Null
$ws.Range("B104:F104").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

$wb.Save()
